$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block B2:G11 is shifted down by one row (old row N's values move to row N+1),
# and a brand-new row of computed values is inserted at row 2.

$cols = @("B", "C", "D", "E", "F", "G")

# Snapshot the current scalar values for rows 2 through 10 (these will move to rows 3 through 11).
$oldValues = @{}
for ($r = 2; $r -le 10; $r++) {
    foreach ($col in $cols) {
        $oldValues["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# Shift rows 10->11, 9->10, ..., 2->3 (process from bottom up so we don't overwrite source data).
for ($r = 10; $r -ge 2; $r--) {
    foreach ($col in $cols) {
        $ws.Range("$col$($r+1)").Value2 = $oldValues["$col$r"]
    }
}

# Write the new first row of data (row 2)
$ws.Range("B2").Value2 = -0.02314597604078636
$ws.Range("C2").Value2 = 0.3579920056255013
$ws.Range("D2").Value2 = 0.1782699060034266
$ws.Range("E2").Value2 = 0.4222202103209018
$ws.Range("F2").Value2 = 0.4363822494547141
$ws.Range("G2").Value2 = 15
